# "Print tensões em Volts + Simulações ANAFAS"
# Insert a new header row above the existing Z0 bus-impedance matrix data
# and label each of the 11 columns "Coluna 1".."Coluna 11".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push all existing data down one row, keeping formatting/links intact.
$ws.Rows(1).Insert()

# Write the new header labels into the freshly inserted row 1.
$headers = @("Coluna 1","Coluna 2","Coluna 3","Coluna 4","Coluna 5","Coluna 6","Coluna 7","Coluna 8","Coluna 9","Coluna 10","Coluna 11")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Match the saved selection state (header row selected).
$ws.Range("A1:K1").Select()
